$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G5").Value = 1.67
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.3
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("X5").Value = 8
$ws.Range("Z5").Value = 13
$ws.Range("AG5").Value = 13
$ws.Range("AO5").Value = 9
$ws.Range("AQ5").Value = 29
$ws.Range("AZ5").Value = 81
$ws.Range("BA5").Value = 101
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 1.88
$ws.Range("BD9").Value = 151
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.3
$ws.Range("Q15").Value = 2.03
$ws.Range("R15").Value = 1.78
$ws.Range("M16").Value = 1.07
$ws.Range("O16").Value = 1.36
$ws.Range("M17").Value = 1.06
$ws.Range("O17").Value = 1.33
$ws.Range("Q17").Value = 2.08
$ws.Range("R17").Value = 1.73
$ws.Range("G24").Value = 2.05
$ws.Range("H24").Value = 3.25
$ws.Range("I24").Value = 3.9
$ws.Range("L24").Value = 4.33
$ws.Range("M24").Value = 1.07
$ws.Range("N24").Value = 9
$ws.Range("Q24").Value = 2.05
$ws.Range("R24").Value = 1.75
$ws.Range("W24").Value = 7
$ws.Range("Z24").Value = 17
$ws.Range("AC24").Value = 9
$ws.Range("AK24").Value = 34
$ws.Range("BA24").Value = 101
$ws.Range("G33").Value = 2.92
$ws.Range("H33").Value = 3.55
$ws.Range("I33").Value = 2.1
$ws.Range("J33").Value = 3.45
$ws.Range("L33").Value = 2.67
$ws.Range("Q33").Value = 1.75
$ws.Range("R33").Value = 2.02
$ws.Range("V33").Value = 2.12
$ws.Range("W33").Value = 10.75
$ws.Range("X33").Value = 16.5
$ws.Range("Y33").Value = 10.75
$ws.Range("Z33").Value = 35
$ws.Range("AA33").Value = 24
$ws.Range("AB33").Value = 29
$ws.Range("AD33").Value = 7.1
$ws.Range("AG33").Value = 8.75
$ws.Range("AH33").Value = 11
$ws.Range("AI33").Value = 8.75
$ws.Range("AJ33").Value = 19.5
$ws.Range("AK33").Value = 16
$ws.Range("AN33").Value = 5
$ws.Range("AO33").Value = 15.5
$ws.Range("AP33").Value = 22
$ws.Range("AQ33").Value = 70
$ws.Range("AR33").Value = 100
$ws.Range("AS33").Value = 250
$ws.Range("AU33").Value = 7
$ws.Range("AV33").Value = 55
$ws.Range("AW33").Value = 4.15
$ws.Range("AX33").Value = 10.75
$ws.Range("AY33").Value = 18
$ws.Range("AZ33").Value = 40
$ws.Range("H37").Value = 3.2
$ws.Range("I37").Value = 3.05
$ws.Range("L37").Value = 3.6
$ws.Range("P37").Value = 2.65
$ws.Range("U37").Value = 1.85
$ws.Range("W37").Value = 6.8
$ws.Range("Z37").Value = 22
$ws.Range("AA37").Value = 19.5
$ws.Range("AD37").Value = 6.2
$ws.Range("AE37").Value = 16
$ws.Range("AG37").Value = 8.25
$ws.Range("AP37").Value = 21
$ws.Range("AT37").Value = 2.4
$ws.Range("AU37").Value = 7.4
$ws.Range("AW37").Value = 4.8
$ws.Range("AX37").Value = 17
$ws.Range("BB37").Value = 400
